# "User can follow sites, Ping : progress"
#
# lifts sheet:
#   - header row (A1:F1) gets a new bold+underline look
#   - lift #0 (row 3) is reset back to an unclaimed/no-note state
#   - the two claimed lift rows (id 1 and id 2) are removed
# follows sheet:
#   - becomes a real "follow a site" list: more site names, plus a
#     Name/Users... header row and a sample "S" site with followers

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: lifts
# ---------------------------------------------------------------
$lifts = $wb.Worksheets.Item("lifts")

# Give the header row its own (bold + underline) style.
$lifts.Range("A1:F1").Font.Bold = $true
$lifts.Range("A1:F1").Font.Underline = $true
$lifts.Range("A1:F1").HorizontalAlignment = -4108
$lifts.Range("A1:F1").VerticalAlignment = -4108

# Lift id=0 (row 3) goes back to NONE / no note, and loses its user.
$lifts.Range("B3").Value = "NONE"
$lifts.Range("E3").Value = "Note"
$lifts.Range("F3").Value = ""

# Lift id=1 and id=2 (rows 4 and 5) are gone.
$lifts.Rows("4:5").Delete()

# ---------------------------------------------------------------
# Sheet: follows
# ---------------------------------------------------------------
$follows = $wb.Worksheets.Item("follows")

# New header: A keeps the bold+underline "title" style, B:E are the
# plain bold header style already used elsewhere in the sheet.
$follows.Range("A1").Value = "Name"
$follows.Range("A1").Font.Bold = $true
$follows.Range("A1").Font.Underline = $true
$follows.Range("A1").HorizontalAlignment = -4108
$follows.Range("A1").VerticalAlignment = -4108

$follows.Range("B1").Value = "Users..."
$follows.Range("C1").Value = "..."
$follows.Range("D1").Value = ".."
$follows.Range("E1").Value = "."
$follows.Range("B1:E1").Font.Bold = $true
$follows.Range("B1:E1").HorizontalAlignment = -4108
$follows.Range("B1:E1").VerticalAlignment = -4108

# Expand the followable-site list (same bold/centered look as the
# existing A2:A3 site names).
$follows.Range("A2").Value = "K3"
$follows.Range("A3").Value = "036G"
$follows.Range("A4").Value = "36G"
$follows.Range("A5").Value = "K4"
$follows.Range("A6").Value = "046P"
$follows.Range("A7").Value = "046G"
$follows.Range("A8").Value = "46G"
$follows.Range("A9").Value = "047G"
$follows.Range("A10").Value = "47G"
$follows.Range("A11").Value = "Pääkeittiö"
$follows.Range("A12").Value = "S"
$follows.Range("A2:A12").Font.Bold = $true
$follows.Range("A2:A12").HorizontalAlignment = -4108
$follows.Range("A2:A12").VerticalAlignment = -4108

# Example row: site "S" already has followers Eemeli and Akseli
# (plain, centered style - like the rest of the non-header cells).
$follows.Range("B12").Value = "Eemeli"
$follows.Range("C12").Value = "Akseli"
$follows.Range("B12:C12").Font.Bold = $false
$follows.Range("B12:C12").HorizontalAlignment = -4108
$follows.Range("B12:C12").VerticalAlignment = -4108

# ---------------------------------------------------------------
# Selection / view bookkeeping (matches what the author last looked at)
# ---------------------------------------------------------------
[void]$follows.Range("F33").Select()
[void]$lifts.Activate()
[void]$lifts.Range("C29").Select()
